# PROS-7407 brand exclude removed
#
# The "brand_name" exclusion list value on the Exclude sheet
# (cells D4, D7, D10, D13) previously held the comma-separated list of
# brands to exclude:
#   "A2 White Milk,Dairy Farmers White Milk,Pauls White Milk,Other Dairy"
# It is cleared back down to an empty list (just the separator ",").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exclude")

$oldValue = "A2 White Milk,Dairy Farmers White Milk,Pauls White Milk,Other Dairy"
$newValue = ","

$cells = @("D4", "D7", "D10", "D13")
foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}

# The active selection on the Exclude sheet moved from D13 to D15.
[void]$ws.Range("D15").Select()

Write-Output "brand_name exclude list cleared"
